# Update column E ("estado") values on the active worksheet to reflect a new
# set of objective-function trial results ("Trabaja" vs "Nada").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "estado" value becomes "Trabaja"
$trabajaRanges = "E2:E6", "E8:E9", "E35", "E39", "E48:E53", "E55", "E80", "E82", "E85", "E94:E98", "E100:E101", "E128", "E140:E144", "E146", "E173", "E175", "E186:E191", "E193", "E218", "E220", "E232:E236", "E238:E239", "E264", "E267", "E269", "E278:E284", "E312", "E315", "E332", "E362:E369"

# Rows whose "estado" value becomes "Nada"
$nadaRanges = "E34", "E40", "E42", "E44:E47", "E75", "E86", "E88", "E90:E93", "E121", "E133", "E135:E139", "E179", "E181:E185", "E213", "E224", "E226", "E228:E231", "E259", "E270", "E272", "E274:E277", "E310", "E316", "E318", "E320:E331", "E337"

foreach ($addr in $trabajaRanges) {
    $ws.Range($addr).Value = "Trabaja"
}

foreach ($addr in $nadaRanges) {
    $ws.Range($addr).Value = "Nada"
}
